$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final player roster table (header + 17 data rows), replacing the
# previous 18-data-row table (the "Ty Jerome" row is dropped and the
# remaining rows are reshuffled/regrouped).
$data = @(
  @("Oyuncu Adı", "Pozisyon", "Takım"),
  @("Donovan Mitchell", "PG,SG", "Cleveland Cavaliers"),
  @("Malik Beasley", "SG,SF", "Detroit Pistons"),
  @("Kelly Oubre Jr.", "SG,SF", "Philadelphia 76ers"),
  @("Josh Hart", "SG,SF,PF", "New York Knicks"),
  @("Jaden McDaniels", "SF,PF", "Minnesota Timberwolves"),
  @("Andrew Wiggins", "SF,PF", "Golden State Warriors"),
  @("Guerschon Yabusele", "PF,C", "Philadelphia 76ers"),
  @("Michael Porter Jr.", "SF,PF", "Denver Nuggets"),
  @("Domantas Sabonis", "C", "Sacramento Kings"),
  @("Victor Wembanyama", "C", "San Antonio Spurs"),
  @("Kel'el Ware", "PF,C", "Miami Heat"),
  @("Kristaps Porzingis", "PF,C", "Boston Celtics"),
  @("Dyson Daniels", "PG,SG,SF", "Atlanta Hawks"),
  @("De'Andre Hunter", "SF,PF", "Atlanta Hawks"),
  @("Alperen Sengün", "C", "Houston Rockets"),
  @("Cam Thomas", "SG,SF", "Brooklyn Nets"),
  @("Donte DiVincenzo", "PG,SG,SF", "Minnesota Timberwolves")
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $i + 1
  $ws.Cells.Item($row, 1).Value = $data[$i][0]
  $ws.Cells.Item($row, 2).Value = $data[$i][1]
  $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# The table shrank by one row (19 -> 18); clear out the old trailing row.
$ws.Range("A19:C19").Clear()
